# Updates cryptos list (automated refresh of Price/Volume(1h) columns,
# plus a few coin rows that got re-ranked and swapped places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: column D (Price) holds text that looks numeric (e.g. "0.900", "99.40",
# "25.826.80" with thousands separators). Setting .Value directly would let
# Excel auto-convert/round these as numbers, so we force text format first,
# assign the value, then restore the default "Normal" style so no stray
# number-format style is left attached to the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.826.80'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.634.59'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.17'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.89'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0781'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.25'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.635.53'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.859.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.559'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₃0768'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.825.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("E21").Value = '  +1.64%  '
$ws.Range("E22").Value = '  +1.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.15'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.55%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("E25").Value = '  -2.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '139.56'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.18%  '
$ws.Range("E27").Value = '  -3.32%  '
$ws.Range("E28").Value = '  +1.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.51'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.69%  '
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0495'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.32'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.31%  '
$ws.Range("E33").Value = '  +1.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.57'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("E35").Value = '  +0.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.900'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("E37").Value = '  +0.51%  '
$ws.Range("E38").Value = '  +0.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.117.64'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0157'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.46%  '
$ws.Range("E41").Value = '  +0.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.55'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.40'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.799'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0₆0108'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.52'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("B47").Value = 'SynthetixNetwork'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +13.63%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.421'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.14%  '
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.61'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.05%  '
